# Insert a new slide ("Test" / "Aantal" / "Feedback") right before the
# final slide of the deck, using the same "Title and Content" layout
# (ppLayoutText = 2, i.e. slideLayout2.xml / "Titel en object") that the
# existing slides in this deck use.
#
# Before: ... , 261, 260 (slide index 6 == last slide, id 260)
# After:  ... , 261, 262(new "Test" slide), 260
$p = $ppt.ActivePresentation
$s = $p.Slides.Add(6, 2)

# Match the Dutch placeholder names used elsewhere in this deck.
$s.Shapes.Item(1).Name = "Titel 1"
$s.Shapes.Item(2).Name = "Tijdelijke aanduiding voor inhoud 2"

# Title placeholder.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Test"

# Content placeholder: two paragraphs, "Aantal" and "Feedback".
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Aantal`rFeedback"
